$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Tema do trabalho do Grupo", $false, $false, $false, $false, $false,
              $true, 1, $false,
              "Kliniek " + [char]0x2013 + " Sistema de Gest" + [char]0x00E3 + "o de Processos Clinicos",
              2)
